$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "DO"
$ws.Range("D2").Value = "M1"
$ws.Range("F2").Value = "M3"
$ws.Range("H2").Value = "M1"
$ws.Range("K2").Value = "M1"
$ws.Range("L2").Value = "M1"
$ws.Range("O2").Value = "M3"
$ws.Range("R2").Value = "M1"
$ws.Range("S2").Value = "DO"
$ws.Range("T2").Value = "M1"
$ws.Range("U2").Value = "M3"
$ws.Range("X2").Value = "M3"
$ws.Range("Y2").Value = "DO"
$ws.Range("AA2").Value = "M1"
$ws.Range("C3").Value = "A1"
$ws.Range("D3").Value = "DO"
$ws.Range("E3").Value = "M1"
$ws.Range("F3").Value = "A2"
$ws.Range("I3").Value = "DO"
$ws.Range("K3").Value = "A1"
$ws.Range("L3").Value = "M1"
$ws.Range("M3").Value = "A2"
$ws.Range("O3").Value = "A1"
$ws.Range("P3").Value = "M3"
$ws.Range("Q3").Value = "A1"
$ws.Range("R3").Value = "A1"
$ws.Range("S3").Value = "M1"
$ws.Range("T3").Value = "DO"
$ws.Range("U3").Value = "M3"
$ws.Range("V3").Value = "A1"
$ws.Range("W3").Value = "A1"
$ws.Range("X3").Value = "DO"
$ws.Range("Z3").Value = "A1"
$ws.Range("AA3").Value = "A1"
$ws.Range("AB3").Value = "M1"
$ws.Range("AC3").Value = "M1"
$ws.Range("B4").Value = "DO"
$ws.Range("C4").Value = "M1"
$ws.Range("D4").Value = "A1"
$ws.Range("E4").Value = "A1"
$ws.Range("F4").Value = "M1"
$ws.Range("H4").Value = "M3"
$ws.Range("I4").Value = "DO"
$ws.Range("K4").Value = "M1"
$ws.Range("L4").Value = "A1"
$ws.Range("M4").Value = "M1"
$ws.Range("O4").Value = "M1"
$ws.Range("Q4").Value = "M1"
$ws.Range("S4").Value = "A1"
$ws.Range("T4").Value = "M1"
$ws.Range("V4").Value = "M3"
$ws.Range("W4").Value = "M1"
$ws.Range("X4").Value = "A1"
$ws.Range("Z4").Value = "DO"
$ws.Range("AA4").Value = "M3"
$ws.Range("AB4").Value = "A1"
$ws.Range("AC4").Value = "A1"
$ws.Range("B5").Value = "M2"
$ws.Range("C5").Value = "DO"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = "A1"
$ws.Range("F5").Value = "A2"
$ws.Range("H5").Value = "M2"
$ws.Range("I5").Value = "M1"
$ws.Range("J5").Value = "A2"
$ws.Range("K5").Value = "M3"
$ws.Range("M5").Value = "A1"
$ws.Range("N5").Value = "M2"
$ws.Range("O5").Value = "DO"
$ws.Range("P5").Value = "DO"
$ws.Range("Q5").Value = "M2"
$ws.Range("R5").Value = "M3"
$ws.Range("S5").Value = "M2"
$ws.Range("T5").Value = "M3"
$ws.Range("U5").Value = "M2"
$ws.Range("Y5").Value = "M1"
$ws.Range("Z5").Value = "M1"
$ws.Range("AA5").Value = "DO"
$ws.Range("AB5").Value = "M2"
$ws.Range("AC5").Value = "M2"
$ws.Range("B6").Value = "M1"
$ws.Range("E6").Value = "M1"
$ws.Range("F6").Value = "M2"
$ws.Range("H6").Value = "DO"
$ws.Range("I6").Value = "M2"
$ws.Range("J6").Value = "M1"
$ws.Range("K6").Value = "A1"
$ws.Range("L6").Value = "M2"
$ws.Range("M6").Value = "A1"
$ws.Range("N6").Value = "M2"
$ws.Range("O6").Value = "DO"
$ws.Range("P6").Value = "M1"
$ws.Range("Q6").Value = "DO"
$ws.Range("R6").Value = "M2"
$ws.Range("S6").Value = "M2"
$ws.Range("T6").Value = "A1"
$ws.Range("V6").Value = "M2"
$ws.Range("AA6").Value = "M1"
$ws.Range("AB6").Value = "M3"
$ws.Range("AC6").Value = "A1"
$ws.Range("F7").Value = "A1"
$ws.Range("G7").Value = "M3"
$ws.Range("J7").Value = "A1"
$ws.Range("K7").Value = "DO"
$ws.Range("T7").Value = "M3"
$ws.Range("U7").Value = "A1"
$ws.Range("V7").Value = "A1"
$ws.Range("AA7").Value = "A1"
$ws.Range("AC7").Value = "M3"
$ws.Range("B8").Value = "DO"
$ws.Range("C8").Value = "M3"
$ws.Range("D8").Value = "A2"
$ws.Range("E8").Value = "A1"
$ws.Range("F8").Value = "M3"
$ws.Range("G8").Value = "A2"
$ws.Range("I8").Value = "M1"
$ws.Range("K8").Value = "A2"
$ws.Range("L8").Value = "M2"
$ws.Range("M8").Value = "M2"
$ws.Range("N8").Value = "M1"
$ws.Range("Q8").Value = "A1"
$ws.Range("R8").Value = "DO"
$ws.Range("S8").Value = "A2"
$ws.Range("T8").Value = "A2"
$ws.Range("U8").Value = "M2"
$ws.Range("V8").Value = "A1"
$ws.Range("W8").Value = "A2"
$ws.Range("X8").Value = "DO"
$ws.Range("Y8").Value = "M2"
$ws.Range("Z8").Value = "A1"
$ws.Range("AA8").Value = "A2"
$ws.Range("AC8").Value = "M2"
$ws.Range("B9").Value = "M1"
$ws.Range("C9").Value = "M2"
$ws.Range("D9").Value = "M1"
$ws.Range("E9").Value = "M1"
$ws.Range("I9").Value = "A2"
$ws.Range("J9").Value = "DO"
$ws.Range("K9").Value = "M2"
$ws.Range("L9").Value = "M2"
$ws.Range("M9").Value = "A1"
$ws.Range("N9").Value = "A1"
$ws.Range("O9").Value = "M2"
$ws.Range("P9").Value = "M1"
$ws.Range("Q9").Value = "M2"
$ws.Range("R9").Value = "A2"
$ws.Range("S9").Value = "M2"
$ws.Range("T9").Value = "A2"
$ws.Range("U9").Value = "A1"
$ws.Range("V9").Value = "DO"
$ws.Range("W9").Value = "A2"
$ws.Range("X9").Value = "M2"
$ws.Range("Y9").Value = "M1"
$ws.Range("Z9").Value = "DO"
$ws.Range("AA9").Value = "M2"
$ws.Range("AB9").Value = "M2"
$ws.Range("AC9").Value = "A1"
$ws.Range("B10").Value = "A2"
$ws.Range("C10").Value = "A1"
$ws.Range("E10").Value = "M2"
$ws.Range("F10").Value = "A2"
$ws.Range("G10").Value = "M1"
$ws.Range("J10").Value = "M1"
$ws.Range("K10").Value = "A1"
$ws.Range("L10").Value = "A2"
$ws.Range("N10").Value = "M2"
$ws.Range("O10").Value = "A2"
$ws.Range("P10").Value = "A1"
$ws.Range("Q10").Value = "M2"
$ws.Range("R10").Value = "M2"
$ws.Range("S10").Value = "DO"
$ws.Range("U10").Value = "M2"
$ws.Range("W10").Value = "M2"
$ws.Range("X10").Value = "A1"
$ws.Range("Y10").Value = "A2"
$ws.Range("Z10").Value = "M1"
$ws.Range("AA10").Value = "M1"
